# Removing less than USD 5 price from extrapolation calibration because it is just a noise
# Updates recalculated values in columns D:H for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = @{ D = 112883.7391489402;  E = 0.004763057682040931;   F = 0.1722540257567458; G = -1.268319468313189;  H = 11.26768481216366 }
    4  = @{ D = 113725.512408094;   E = -0.0005611441198597482; F = 0.1908204657038853; G = -1.509263841052931;  H = 12.59163667277719 }
    5  = @{ D = 114455.0058600235;  E = -0.003820815790867739;  F = 0.1951824710473767; G = -0.8102660972289101; H = 7.585638409631578 }
    6  = @{ D = 115060.2364792997;  E = -0.01300777817242779;   F = 0.2293912338860257; G = -1.386852238559575;  H = 11.07136796359288 }
    8  = @{ D = 116907.7111220563;  E = -0.03258664512340192;   F = 0.2086998108458011; G = -1.16086155669494;   H = 7.817787217133953 }
    10 = @{ D = 119967.0285379993;  E = -0.1056033928387945;    F = 0.4405374640767752; G = -1.88591738007301;   H = 9.525320965214744 }
    11 = @{ D = 121648.2729414744;  E = -0.1824791195662304;    F = 0.7530824280261705; G = -2.50555008259193;   H = 11.87187476787148 }
    13 = @{ D = 112872.0668740137;  E = 0.05367345592475129;    F = 0.1358023128035912; G = -0.8293083886466766; H = 6.780960180530164 }
    14 = @{ D = 112889.342296171;   E = 0.04929407705994201;    F = 0.1380951289343995; G = -0.7812796165832036; H = 6.837270817202463 }
    16 = @{ D = 112904.9703018227;  E = 0.04571232644541313;    F = 0.1387826665003753; G = -0.7267543967038677; H = 7.492427315641736 }
    17 = @{ D = 112821.5129246753;  E = 0.03663106619014705;    F = 0.1442794364729697; G = -1.041845963693458;  H = 6.132679052189135 }
    18 = @{ D = 112958.9584880201;  E = 0.02592122804018062;    F = 0.1532453225384673; G = -0.6947024462290753; H = 7.118589828545098 }
    19 = @{ D = 112969.7129666251;  E = 0.01429496541325808;    F = 0.1642713640594899; G = -0.6561811003281649; H = 6.689541172995616 }
    20 = @{ D = 113780.3037758657;  E = 0.01636482599261388;    F = 0.2377040385869526; G = -4.041282978400303;  H = 42.48128837360665 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
}
